# Applies transition-matrix probability updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "B2" = 0.2352941176470588
    "C2" = 0.5294117647058824
    "P2" = 0.1176470588235294
    "S2" = 0.1176470588235294
    "C3" = 0.1
    "P3" = 0.5
    "S3" = 0.4
    "P4" = 1.0
    "D6" = 0.07142857142857142
    "F6" = 0.07142857142857142
    "J6" = 0.3571428571428572
    "Q6" = 0.1428571428571428
    "R6" = 0.07142857142857142
    "S6" = 0.2857142857142857
    "J7" = 0.1666666666666667
    "Q7" = 0.08333333333333333
    "S7" = 0.75
    "B8" = 0.05263157894736842
    "F8" = 0.05263157894736842
    "J8" = 0.07894736842105263
    "Q8" = 0.02631578947368421
    "R8" = 0.07894736842105263
    "S8" = 0.7105263157894737
    "B9" = 0.04545454545454546
    "F9" = 0.09090909090909091
    "J9" = 0.1363636363636364
    "Q9" = 0.1363636363636364
    "R9" = 0.04545454545454546
    "S9" = 0.5454545454545454
    "B10" = 0.1176470588235294
    "D10" = 0.01176470588235294
    "F10" = 0.04705882352941176
    "J10" = 0.1176470588235294
    "O10" = 0.01176470588235294
    "Q10" = 0.2235294117647059
    "R10" = 0.09411764705882353
    "S10" = 0.3764705882352941
    "G11" = 0.08695652173913043
    "J11" = 0.1739130434782609
    "K11" = 0.1304347826086956
    "L11" = 0.6086956521739131
    "G12" = 0.7142857142857143
    "J12" = 0.1428571428571428
    "S12" = 0.1428571428571428
    "G13" = 0.5
    "S13" = 0.5
    "J14" = 1.0
    "F15" = 0.06666666666666667
    "H15" = 0.1333333333333333
    "I15" = 0.06666666666666667
    "J15" = 0.2
    "K15" = 0.1333333333333333
    "O15" = 0.06666666666666667
    "S15" = 0.3333333333333333
    "H16" = 0.5555555555555556
    "I16" = 0.1111111111111111
    "J16" = 0.2222222222222222
    "S16" = 0.1111111111111111
    "F17" = 0.03846153846153846
    "H17" = 0.1538461538461539
    "I17" = 0.03846153846153846
    "J17" = 0.3846153846153846
    "M17" = 0.03846153846153846
    "N17" = 0.03846153846153846
    "O17" = 0.1153846153846154
    "S17" = 0.1923076923076923
    "H18" = 0.3076923076923077
    "I18" = 0.2307692307692308
    "J18" = 0.1538461538461539
    "K18" = 0.1538461538461539
    "S18" = 0.1538461538461539
    "F19" = 0.01652892561983471
    "H19" = 0.1900826446280992
    "I19" = 0.1322314049586777
    "J19" = 0.3223140495867768
    "K19" = 0.1322314049586777
    "M19" = 0.008264462809917356
    "O19" = 0.05785123966942149
    "S19" = 0.140495867768595
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

